# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column D) for the
# 49ed7b88-ffcd-4894-879d-8a6b41754fa8 file row (row 5) on both the
# zh-cn and de-de localization-status sheets, reflecting a new handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-01 08:54:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-01 08:54:57"
